$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 55.08
$ws.Range("E2").Value = 65.6
$ws.Range("F2").Value = 17.44
$ws.Range("K2").Value = 58.2
$ws.Range("N2").Value = 54.02451352198364

# Row 3
$ws.Range("D3").Value = 29.64
$ws.Range("E3").Value = 62.3
$ws.Range("F3").Value = 15.9
$ws.Range("K3").Value = 56.4
$ws.Range("N3").Value = 54.02451352198364

# Row 4
$ws.Range("D4").Value = 28.51
$ws.Range("E4").Value = 65.4
$ws.Range("F4").Value = 27.21
$ws.Range("K4").Value = 56.4
$ws.Range("N4").Value = 54.02451352198364

# Row 5
$ws.Range("D5").Value = 307.57
$ws.Range("E5").Value = 52.5
$ws.Range("F5").Value = 1.44
$ws.Range("K5").Value = 56.4
$ws.Range("N5").Value = 54.02451352198364
